# Auto-generated script to update Chocobo Profits workbook market data
# Applies per-cell numeric updates to sheets ALC, ARM, BSM, CRP, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: ALC (18 cell updates) ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("H96").Value = 862.5
$ws.Range("I96").Value = 275
$ws.Range("J96").Value = 1450
$ws.Range("K96").Value = 825
$ws.Range("L96").Value = 4350
$ws.Range("M96").Value = 548
$ws.Range("N96").Value = -7096
$ws.Range("H132").Value = 19310252
$ws.Range("I132").Value = 21363592
$ws.Range("J132").Value = 8868
$ws.Range("K132").Value = 64090776
$ws.Range("L132").Value = 26604
$ws.Range("M132").Value = -64088246
$ws.Range("N132").Value = -31664
$ws.Range("H137").Value = 2839.138
$ws.Range("I137").Value = 1853.4
$ws.Range("K137").Value = 5560.200000000001
$ws.Range("M137").Value = -3010.200000000001

# ---- Sheet 2: ARM (59 cell updates) ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 6760.2656
$ws.Range("I32").Value = 4121.4
$ws.Range("K32").Value = 4121.4
$ws.Range("M32").Value = -3834.4
$ws.Range("H45").Value = 1670.2222
$ws.Range("I45").Value = 1317.7142
$ws.Range("J45").Value = 1894.5454
$ws.Range("K45").Value = 1317.7142
$ws.Range("L45").Value = 1894.5454
$ws.Range("M45").Value = -940.7141999999999
$ws.Range("N45").Value = -2648.5454
$ws.Range("H61").Value = 2767.4443
$ws.Range("I61").Value = 2140
$ws.Range("K61").Value = 2140
$ws.Range("M61").Value = -1928
$ws.Range("H74").Value = 2138.8276
$ws.Range("I74").Value = 1567.3043
$ws.Range("J74").Value = 4329.6665
$ws.Range("K74").Value = 1567.3043
$ws.Range("L74").Value = 4329.6665
$ws.Range("M74").Value = -693.3043
$ws.Range("N74").Value = -6077.6665
$ws.Range("H77").Value = 2138.8276
$ws.Range("I77").Value = 1567.3043
$ws.Range("J77").Value = 4329.6665
$ws.Range("K77").Value = 7836.5215
$ws.Range("L77").Value = 21648.3325
$ws.Range("M77").Value = -3468.5215
$ws.Range("N77").Value = -30384.3325
$ws.Range("H110").Value = 867.3125
$ws.Range("I110").Value = 847.9583
$ws.Range("J110").Value = 925.375
$ws.Range("K110").Value = 847.9583
$ws.Range("L110").Value = 925.375
$ws.Range("M110").Value = 1197.0417
$ws.Range("N110").Value = -5015.375
$ws.Range("H115").Value = 29880
$ws.Range("J115").Value = 29880
$ws.Range("L115").Value = 29880
$ws.Range("N115").Value = -33014
$ws.Range("H122").Value = 2108.6
$ws.Range("I122").Value = 1212.8235
$ws.Range("J122").Value = 3280
$ws.Range("K122").Value = 3638.4705
$ws.Range("L122").Value = 9840
$ws.Range("M122").Value = -1188.4705
$ws.Range("N122").Value = -14740
$ws.Range("H132").Value = 2732.6135
$ws.Range("I132").Value = 1972.1072
$ws.Range("K132").Value = 5916.321599999999
$ws.Range("M132").Value = -3386.321599999999
$ws.Range("H136").Value = 2767.4443
$ws.Range("I136").Value = 2140
$ws.Range("K136").Value = 6420
$ws.Range("M136").Value = -3870
$ws.Range("H137").Value = 50548
$ws.Range("J137").Value = 50548
$ws.Range("L137").Value = 50548
$ws.Range("N137").Value = -60748

# ---- Sheet 3: BSM (15 cell updates) ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("H27").Value = 28000
$ws.Range("J27").Value = 28000
$ws.Range("L27").Value = 28000
$ws.Range("N27").Value = -28384
$ws.Range("H107").Value = 1060.7142
$ws.Range("I107").Value = 1020.8333
$ws.Range("K107").Value = 1020.8333
$ws.Range("M107").Value = 899.1667
$ws.Range("H134").Value = 3837.8064
$ws.Range("I134").Value = 2234.1177
$ws.Range("J134").Value = 5785.143
$ws.Range("K134").Value = 6702.353099999999
$ws.Range("L134").Value = 17355.429
$ws.Range("M134").Value = -4167.353099999999
$ws.Range("N134").Value = -22425.429

# ---- Sheet 4: CRP (42 cell updates) ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 10872060
$ws.Range("I31").Value = 1008.15625
$ws.Range("J31").Value = 35720180
$ws.Range("K31").Value = 1008.15625
$ws.Range("L31").Value = 35720180
$ws.Range("M31").Value = -713.15625
$ws.Range("N31").Value = -35720770
$ws.Range("H34").Value = 10872060
$ws.Range("I34").Value = 1008.15625
$ws.Range("J34").Value = 35720180
$ws.Range("K34").Value = 1008.15625
$ws.Range("L34").Value = 35720180
$ws.Range("M34").Value = -806.15625
$ws.Range("N34").Value = -35720584
$ws.Range("H58").Value = 1871.3611
$ws.Range("I58").Value = 1641.0968
$ws.Range("J58").Value = 3299
$ws.Range("K58").Value = 1641.0968
$ws.Range("L58").Value = 3299
$ws.Range("M58").Value = -1438.0968
$ws.Range("N58").Value = -3705
$ws.Range("H132").Value = 4132.75
$ws.Range("I132").Value = 3988.0715
$ws.Range("J132").Value = 4224.8184
$ws.Range("K132").Value = 11964.2145
$ws.Range("L132").Value = 12674.4552
$ws.Range("M132").Value = -9434.2145
$ws.Range("N132").Value = -17734.4552
$ws.Range("H134").Value = 8471.412
$ws.Range("I134").Value = 13731
$ws.Range("J134").Value = 3796.2222
$ws.Range("K134").Value = 41193
$ws.Range("L134").Value = 11388.6666
$ws.Range("M134").Value = -38658
$ws.Range("N134").Value = -16458.6666
$ws.Range("H136").Value = 1871.3611
$ws.Range("I136").Value = 1641.0968
$ws.Range("J136").Value = 3299
$ws.Range("K136").Value = 4923.2904
$ws.Range("L136").Value = 9897
$ws.Range("M136").Value = -2373.2904
$ws.Range("N136").Value = -14997

# ---- Sheet 6: GSM (32 cell updates) ----
$ws = $wb.Worksheets.Item(6)
$ws.Range("H113").Value = 1301.091
$ws.Range("I113").Value = 1269.125
$ws.Range("J113").Value = 1386.3334
$ws.Range("K113").Value = 1269.125
$ws.Range("L113").Value = 1386.3334
$ws.Range("M113").Value = 900.875
$ws.Range("N113").Value = -5726.3334
$ws.Range("H122").Value = 2819.1875
$ws.Range("I122").Value = 1920.7
$ws.Range("J122").Value = 4316.6665
$ws.Range("K122").Value = 5762.1
$ws.Range("L122").Value = 12949.9995
$ws.Range("M122").Value = -3312.1
$ws.Range("N122").Value = -17849.9995
$ws.Range("H126").Value = 4014.9102
$ws.Range("I126").Value = 2737.186
$ws.Range("J126").Value = 5209.304
$ws.Range("K126").Value = 8211.558000000001
$ws.Range("L126").Value = 15627.912
$ws.Range("M126").Value = -5741.558000000001
$ws.Range("N126").Value = -20567.912
$ws.Range("H132").Value = 3156.8462
$ws.Range("I132").Value = 1993.3889
$ws.Range("J132").Value = 5774.625
$ws.Range("K132").Value = 5980.1667
$ws.Range("L132").Value = 17323.875
$ws.Range("M132").Value = -3450.1667
$ws.Range("N132").Value = -22383.875
$ws.Range("H137").Value = 63751.6
$ws.Range("J137").Value = 63751.6
$ws.Range("L137").Value = 63751.6
$ws.Range("N137").Value = -73951.60000000001

# ---- Sheet 7: LTW (43 cell updates) ----
$ws = $wb.Worksheets.Item(7)
$ws.Range("H61").Value = 1480.4286
$ws.Range("I61").Value = 1474.3125
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 1474.3125
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -1272.3125
$ws.Range("N61").Value = -1904
$ws.Range("H69").Value = 520000
$ws.Range("J69").Value = 520000
$ws.Range("L69").Value = 520000
$ws.Range("N69").Value = -521622
$ws.Range("H72").Value = 520000
$ws.Range("J72").Value = 520000
$ws.Range("L72").Value = 1560000
$ws.Range("N72").Value = -1568112
$ws.Range("H113").Value = 1480.4286
$ws.Range("I113").Value = 1474.3125
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1474.3125
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 695.6875
$ws.Range("N113").Value = -5840
$ws.Range("H122").Value = 4157.4585
$ws.Range("I122").Value = 2649.625
$ws.Range("J122").Value = 7173.125
$ws.Range("K122").Value = 7948.875
$ws.Range("L122").Value = 21519.375
$ws.Range("M122").Value = -5498.875
$ws.Range("N122").Value = -26419.375
$ws.Range("H132").Value = 3640.6033
$ws.Range("I132").Value = 2637.75
$ws.Range("J132").Value = 5963
$ws.Range("K132").Value = 7913.25
$ws.Range("L132").Value = 17889
$ws.Range("M132").Value = -5383.25
$ws.Range("N132").Value = -22949
$ws.Range("H136").Value = 3686.4614
$ws.Range("I136").Value = 1791.7368
$ws.Range("J136").Value = 8829.286
$ws.Range("K136").Value = 5375.2104
$ws.Range("L136").Value = 26487.858
$ws.Range("M136").Value = -2825.2104
$ws.Range("N136").Value = -31587.858

# ---- Sheet 8: WVR (28 cell updates) ----
$ws = $wb.Worksheets.Item(8)
$ws.Range("H107").Value = 740.75
$ws.Range("I107").Value = 514.3333
$ws.Range("J107").Value = 1148.3
$ws.Range("K107").Value = 1542.9999
$ws.Range("L107").Value = 3444.9
$ws.Range("M107").Value = 377.0001
$ws.Range("N107").Value = -7284.9
$ws.Range("H122").Value = 3884
$ws.Range("I122").Value = 2778.25
$ws.Range("J122").Value = 5358.3335
$ws.Range("K122").Value = 8334.75
$ws.Range("L122").Value = 16075.0005
$ws.Range("M122").Value = -5884.75
$ws.Range("N122").Value = -20975.0005
$ws.Range("H132").Value = 15154232
$ws.Range("I132").Value = 872
$ws.Range("J132").Value = 19611104
$ws.Range("K132").Value = 2616
$ws.Range("L132").Value = 58833312
$ws.Range("M132").Value = -86
$ws.Range("N132").Value = -58838372
$ws.Range("H136").Value = 1597.2759
$ws.Range("I136").Value = 877.1905
$ws.Range("J136").Value = 3487.5
$ws.Range("K136").Value = 2631.5715
$ws.Range("L136").Value = 10462.5
$ws.Range("M136").Value = -81.57150000000001
$ws.Range("N136").Value = -15562.5
